# Append new log rows (207-215) to the logs sheet, matching the
# "Завершить общение" (end-chat) button feature log entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 207; A = 205; B = 206; C = 3; D = "Начал взаимодействие с консультантом!"; E = "09/06/2023 17:50:00" },
    @{ Row = 208; A = 206; B = 207; C = 5; D = "Успешно добавлен в базу!";              E = "09/06/2023 18:11:15" },
    @{ Row = 209; A = 207; B = 208; C = 5; D = "Начал взаимодействие с консультантом!"; E = "09/06/2023 18:11:18" },
    @{ Row = 210; A = 208; B = 209; C = 1; D = "Начал взаимодействие с консультантом!"; E = "11/06/2023 21:33:25" },
    @{ Row = 211; A = 209; B = 210; C = 1; D = "Начал взаимодействие с консультантом!"; E = "11/06/2023 23:11:18" },
    @{ Row = 212; A = 210; B = 211; C = 1; D = "Начал взаимодействие с консультантом!"; E = "11/06/2023 23:17:40" },
    @{ Row = 213; A = 211; B = 212; C = 1; D = "Начал взаимодействие с консультантом!"; E = "11/06/2023 23:24:24" },
    @{ Row = 214; A = 212; B = 213; C = 1; D = "Начал взаимодействие с консультантом!"; E = "11/06/2023 23:26:12" },
    @{ Row = 215; A = 213; B = 214; C = 1; D = "Начал взаимодействие с консультантом!"; E = "11/06/2023 23:34:08" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Insert a fresh row at the bottom of the table (mirrors the row above
    # so the new row inherits the sheet's default formatting baseline).
    $ws.Rows.Item($rowNum).Insert()

    $a = $ws.Cells.Item($rowNum, 1)
    $a.Value = $r.A
    $a.Font.Bold = $true
    $a.HorizontalAlignment = -4108
    $a.VerticalAlignment = -4160
    $a.Borders.LineStyle = 1

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
}
